# "pushing 5-5 from laptop"
# Updates the "Ordered" sheet: adds a "Percentage Events" column (N),
# recomputes the per-column totals row (row 4), and re-categorizes a
# number of T/W/Uncategorized cells after the underlying raw data changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ordered")

# --- New column N: "Percentage Events" header + per-row percentage formulas ---
$ws.Range("N1").Value = "Percentage Events"
$ws.Range("N2").Formula = "=L2/L4"
$ws.Range("N3").Formula = "=L3/L4"

# --- Row 4: replace the old category-letter row with the new per-column totals ---
$ws.Range("B4").Value = 116
$ws.Range("C4").Value = 108
$ws.Range("D4").Value = 133
$ws.Range("E4").Value = 150
$ws.Range("F4").Value = 119
$ws.Range("G4").Value = 125
$ws.Range("H4").Value = 116
$ws.Range("I4").Value = 140
$ws.Range("J4").Value = 196
$ws.Range("K4").Value = 219
$ws.Range("L4").Formula = "=SUM(B4:K4)"

# --- Individual re-categorized cells (sorted T/W/C/Uncategorized blocks shifted) ---
$ws.Range("B6").Value = "C"
$ws.Range("G7").Value = "T"
$ws.Range("K8").Value = "C"
$ws.Range("D10").Value = "T"
$ws.Range("I10").Value = "T"
$ws.Range("E11").Value = "T"
$ws.Range("F11").Value = "T"
$ws.Range("J11").Value = "T"
$ws.Range("H15").Value = "T"
$ws.Range("C17").Value = "T"
$ws.Range("B21").Value = "T"
$ws.Range("K31").Value = "T"

$ws.Range("C106").Value = "Uncategorized"
$ws.Range("H109").Value = "Uncategorized"
$ws.Range("B111").Value = "Uncategorized"
$ws.Range("C112").Value = "W"
$ws.Range("G117").Value = "Uncategorized"
$ws.Range("F119").Value = "Uncategorized"
$ws.Range("B120").Value = "W"
$ws.Range("H120").Value = "W"
$ws.Range("F123").Value = "W"
$ws.Range("D127").Value = "Uncategorized"
$ws.Range("G129").Value = "W"
$ws.Range("I132").Value = "Uncategorized"
$ws.Range("D137").Value = "W"
$ws.Range("I144").Value = "W"
$ws.Range("E147").Value = "Uncategorized"
$ws.Range("E154").Value = "W"
$ws.Range("J191").Value = "Uncategorized"
$ws.Range("J200").Value = "W"
$ws.Range("K203").Value = "Uncategorized"
$ws.Range("K223").Value = "W"

# --- Re-apply the autofilter / sort-range bookkeeping over the now-shifted table ---
$ws.Range("N11").Select()
